$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 62,6
$arr[0,0] = 2025
$arr[0,1] = "CH"
$arr[0,2] = "M2"
$arr[0,3] = "LIM3"
$arr[0,4] = "F"
$arr[0,5] = 98
$arr[1,0] = 2025
$arr[1,1] = "CH"
$arr[1,2] = "M1"
$arr[1,3] = "LIM4"
$arr[1,4] = "M"
$arr[1,5] = 75
$arr[2,0] = 2025
$arr[2,1] = "FT"
$arr[2,2] = "L3"
$arr[2,3] = "LIEEA_AII3"
$arr[2,4] = "F"
$arr[2,5] = 16
$arr[3,0] = 2025
$arr[3,1] = "FT"
$arr[3,2] = "L3"
$arr[3,3] = "LIEEA_AII3"
$arr[3,4] = "M"
$arr[3,5] = 13
$arr[4,0] = 2025
$arr[4,1] = "FT"
$arr[4,2] = "L3"
$arr[4,3] = "LIGE_ERE3"
$arr[4,4] = "F"
$arr[4,5] = 23
$arr[5,0] = 2025
$arr[5,1] = "FT"
$arr[5,2] = "L3"
$arr[5,3] = "LIGE_ERE3"
$arr[5,4] = "M"
$arr[5,5] = 13
$arr[6,0] = 2025
$arr[6,1] = "SV"
$arr[6,2] = "L3"
$arr[6,3] = "BMC3"
$arr[6,4] = "F"
$arr[6,5] = 18
$arr[7,0] = 2025
$arr[7,1] = "SV"
$arr[7,2] = "L3"
$arr[7,3] = "BMC3"
$arr[7,4] = "M"
$arr[7,5] = 2
$arr[8,0] = 2025
$arr[8,1] = "SV"
$arr[8,2] = "L3"
$arr[8,3] = "LIB3"
$arr[8,4] = "F"
$arr[8,5] = 53
$arr[9,0] = 2025
$arr[9,1] = "SV"
$arr[9,2] = "L3"
$arr[9,3] = "LIB3"
$arr[9,4] = "M"
$arr[9,5] = 4
$arr[10,0] = 2025
$arr[10,1] = "SV"
$arr[10,2] = "L3"
$arr[10,3] = "LISVT3"
$arr[10,4] = "F"
$arr[10,5] = 18
$arr[11,0] = 2025
$arr[11,1] = "SV"
$arr[11,2] = "L3"
$arr[11,3] = "LISVT3"
$arr[11,4] = "M"
$arr[11,5] = 2
$arr[12,0] = 2025
$arr[12,1] = "PHYS"
$arr[12,2] = "L2"
$arr[12,3] = "LIM3"
$arr[12,4] = "F"
$arr[12,5] = 12
$arr[13,0] = 2025
$arr[13,1] = "CH"
$arr[13,2] = "L3"
$arr[13,3] = "LIM4"
$arr[13,4] = "M"
$arr[13,5] = 22
$arr[14,0] = 2023
$arr[14,1] = "Informatique"
$arr[14,2] = "M2"
$arr[14,3] = "LIM1"
$arr[14,4] = "F"
$arr[14,5] = 110
$arr[15,0] = 2023
$arr[15,1] = "SV"
$arr[15,2] = "M2"
$arr[15,3] = "LIM2"
$arr[15,4] = "M"
$arr[15,5] = 120
$arr[16,0] = 2023
$arr[16,1] = "CH"
$arr[16,2] = "M2"
$arr[16,3] = "LIM3"
$arr[16,4] = "F"
$arr[16,5] = 98
$arr[17,0] = 2023
$arr[17,1] = "CH"
$arr[17,2] = "M1"
$arr[17,3] = "LIM4"
$arr[17,4] = "M"
$arr[17,5] = 75
$arr[18,0] = 2023
$arr[18,1] = "FT"
$arr[18,2] = "L3"
$arr[18,3] = "LIEEA_AII3"
$arr[18,4] = "F"
$arr[18,5] = 16
$arr[19,0] = 2023
$arr[19,1] = "FT"
$arr[19,2] = "L3"
$arr[19,3] = "LIEEA_AII3"
$arr[19,4] = "F"
$arr[19,5] = 16
$arr[20,0] = 2023
$arr[20,1] = "FT"
$arr[20,2] = "L3"
$arr[20,3] = "LIEEA_AII3"
$arr[20,4] = "M"
$arr[20,5] = 13
$arr[21,0] = 2023
$arr[21,1] = "FT"
$arr[21,2] = "L3"
$arr[21,3] = "LIGE_ERE3"
$arr[21,4] = "F"
$arr[21,5] = 23
$arr[22,0] = 2023
$arr[22,1] = "FT"
$arr[22,2] = "L3"
$arr[22,3] = "LIGE_ERE3"
$arr[22,4] = "M"
$arr[22,5] = 13
$arr[23,0] = 2023
$arr[23,1] = "SV"
$arr[23,2] = "L3"
$arr[23,3] = "BMC3"
$arr[23,4] = "F"
$arr[23,5] = 18
$arr[24,0] = 2023
$arr[24,1] = "SV"
$arr[24,2] = "L3"
$arr[24,3] = "BMC3"
$arr[24,4] = "M"
$arr[24,5] = 2
$arr[25,0] = 2023
$arr[25,1] = "SV"
$arr[25,2] = "L3"
$arr[25,3] = "LIB3"
$arr[25,4] = "F"
$arr[25,5] = 53
$arr[26,0] = 2023
$arr[26,1] = "SV"
$arr[26,2] = "L3"
$arr[26,3] = "LIB3"
$arr[26,4] = "M"
$arr[26,5] = 4
$arr[27,0] = 2023
$arr[27,1] = "SV"
$arr[27,2] = "L3"
$arr[27,3] = "LISVT3"
$arr[27,4] = "F"
$arr[27,5] = 18
$arr[28,0] = 2023
$arr[28,1] = "SV"
$arr[28,2] = "L3"
$arr[28,3] = "LISVT3"
$arr[28,4] = "M"
$arr[28,5] = 2
$arr[29,0] = 2023
$arr[29,1] = "INFO"
$arr[29,2] = "L3"
$arr[29,3] = "IMM3"
$arr[29,4] = "F"
$arr[29,5] = 7
$arr[30,0] = 2023
$arr[30,1] = "INFO"
$arr[30,2] = "L3"
$arr[30,3] = "IMM3"
$arr[30,4] = "M"
$arr[30,5] = 8
$arr[31,0] = 2023
$arr[31,1] = "PHYS"
$arr[31,2] = "L2"
$arr[31,3] = "LIM3"
$arr[31,4] = "F"
$arr[31,5] = 12
$arr[32,0] = 2023
$arr[32,1] = "CH"
$arr[32,2] = "L3"
$arr[32,3] = "LIM4"
$arr[32,4] = "M"
$arr[32,5] = 22
$arr[33,0] = 2023
$arr[33,1] = "PHYS"
$arr[33,2] = "M1"
$arr[33,3] = "LIM5"
$arr[33,4] = "F"
$arr[33,5] = 66
$arr[34,0] = 2023
$arr[34,1] = "CH"
$arr[34,2] = "M1"
$arr[34,3] = "Info"
$arr[34,4] = "F"
$arr[34,5] = 145
$arr[35,0] = 2023
$arr[35,1] = "PHYS"
$arr[35,2] = "M1"
$arr[35,3] = "Info"
$arr[35,4] = "M"
$arr[35,5] = 145
$arr[36,0] = 2023
$arr[36,1] = "Informatique"
$arr[36,2] = "M2"
$arr[36,3] = "LIM1"
$arr[36,4] = "F"
$arr[36,5] = 110
$arr[37,0] = 2023
$arr[37,1] = "SV"
$arr[37,2] = "M2"
$arr[37,3] = "LIM2"
$arr[37,4] = "M"
$arr[37,5] = 120
$arr[38,0] = 2023
$arr[38,1] = "CH"
$arr[38,2] = "M2"
$arr[38,3] = "LIM3"
$arr[38,4] = "F"
$arr[38,5] = 98
$arr[39,0] = 2023
$arr[39,1] = "CH"
$arr[39,2] = "M1"
$arr[39,3] = "LIM4"
$arr[39,4] = "M"
$arr[39,5] = 75
$arr[40,0] = 2022
$arr[40,1] = "SV"
$arr[40,2] = "L3"
$arr[40,3] = "LIB3"
$arr[40,4] = "M"
$arr[40,5] = 4
$arr[41,0] = 2022
$arr[41,1] = "SV"
$arr[41,2] = "L3"
$arr[41,3] = "LISVT3"
$arr[41,4] = "F"
$arr[41,5] = 44
$arr[42,0] = 2022
$arr[42,1] = "SV"
$arr[42,2] = "L3"
$arr[42,3] = "LISVT3"
$arr[42,4] = "M"
$arr[42,5] = 33
$arr[43,0] = 2022
$arr[43,1] = "INFO"
$arr[43,2] = "L3"
$arr[43,3] = "IMM3"
$arr[43,4] = "F"
$arr[43,5] = 22
$arr[44,0] = 2022
$arr[44,1] = "INFO"
$arr[44,2] = "L3"
$arr[44,3] = "IMM3"
$arr[44,4] = "M"
$arr[44,5] = 55
$arr[45,0] = 2022
$arr[45,1] = "PHYS"
$arr[45,2] = "L2"
$arr[45,3] = "LIM3"
$arr[45,4] = "F"
$arr[45,5] = 701
$arr[46,0] = 2022
$arr[46,1] = "CH"
$arr[46,2] = "L3"
$arr[46,3] = "LIM4"
$arr[46,4] = "M"
$arr[46,5] = 54
$arr[47,0] = 2022
$arr[47,1] = "PHYS"
$arr[47,2] = "M1"
$arr[47,3] = "LIM5"
$arr[47,4] = "F"
$arr[47,5] = 12
$arr[48,0] = 2022
$arr[48,1] = "CH"
$arr[48,2] = "M1"
$arr[48,3] = "Info"
$arr[48,4] = "F"
$arr[48,5] = 14
$arr[49,0] = 2022
$arr[49,1] = "PHYS"
$arr[49,2] = "M1"
$arr[49,3] = "Info"
$arr[49,4] = "M"
$arr[49,5] = 42
$arr[50,0] = 2022
$arr[50,1] = "Informatique"
$arr[50,2] = "M2"
$arr[50,3] = "LIM1"
$arr[50,4] = "F"
$arr[50,5] = 51
$arr[51,0] = 2022
$arr[51,1] = "SV"
$arr[51,2] = "M2"
$arr[51,3] = "LIM2"
$arr[51,4] = "M"
$arr[51,5] = 21
$arr[52,0] = 2022
$arr[52,1] = "CH"
$arr[52,2] = "M2"
$arr[52,3] = "LIM3"
$arr[52,4] = "F"
$arr[52,5] = 32
$arr[53,0] = 2022
$arr[53,1] = "CH"
$arr[53,2] = "M1"
$arr[53,3] = "LIM4"
$arr[53,4] = "M"
$arr[53,5] = 12
$arr[54,0] = 2022
$arr[54,1] = "FT"
$arr[54,2] = "L3"
$arr[54,3] = "LIEEA_AII3"
$arr[54,4] = "F"
$arr[54,5] = 15
$arr[55,0] = 2022
$arr[55,1] = "FT"
$arr[55,2] = "L3"
$arr[55,3] = "LIEEA_AII3"
$arr[55,4] = "M"
$arr[55,5] = 18
$arr[56,0] = 2022
$arr[56,1] = "FT"
$arr[56,2] = "L3"
$arr[56,3] = "LIGE_ERE3"
$arr[56,4] = "F"
$arr[56,5] = 23
$arr[57,0] = 2022
$arr[57,1] = "FT"
$arr[57,2] = "L3"
$arr[57,3] = "LIGE_ERE3"
$arr[57,4] = "M"
$arr[57,5] = 13
$arr[58,0] = 2022
$arr[58,1] = "SV"
$arr[58,2] = "L3"
$arr[58,3] = "BMC3"
$arr[58,4] = "F"
$arr[58,5] = 18
$arr[59,0] = 2022
$arr[59,1] = "SV"
$arr[59,2] = "L3"
$arr[59,3] = "BMC3"
$arr[59,4] = "M"
$arr[59,5] = 2
$arr[60,0] = 2022
$arr[60,1] = "SV"
$arr[60,2] = "L3"
$arr[60,3] = "LIB3"
$arr[60,4] = "F"
$arr[60,5] = 53
$arr[61,0] = 2022
$arr[61,1] = "SV"
$arr[61,2] = "L3"
$arr[61,3] = "LIB3"
$arr[61,4] = "M"
$arr[61,5] = 4

$ws.Range("A119:F180").Value() = $arr

$ws.Range("J175").Select()